$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert the two new rows that the final layout needs:
#   - row 9  -> new data row for "Graf" (pushes old row 10.. down by one)
#   - row 17 -> (after the first insert, the old row 16 sits at 17) a second
#               blank spacer row, pushing the old row 16.. down by one more
# This reproduces the exact final row numbering used in the target sheet:
#   9, 11-15, 18-20, 22-23, 26-28, 30
# ---------------------------------------------------------------------------
$ws.Rows("9:9").Insert()
$ws.Rows("17:17").Insert()

# ---------------------------------------------------------------------------
# New data row for Graf (row 9)
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "Graf"
$ws.Range("B9").Value = 4
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 48
$ws.Range("F9").Formula = "=E9*3"

# ---------------------------------------------------------------------------
# F2:F8 - turn the per-row formula into one shared formula group
# ---------------------------------------------------------------------------
$ws.Range("F2:F8").Formula = "=E2*3"

# ---------------------------------------------------------------------------
# Totals rows (formerly rows 10-14, now 11-15) - widen ranges to include the
# new Graf row (row 9)
# ---------------------------------------------------------------------------
$ws.Range("B11").Formula = "=B2+B3+B4+B5+B6+B7+B8+B9"
$ws.Range("C11").Formula = "=AVERAGE(B2:B9)"
$ws.Range("D11").Formula = "=_xlfn.STDEV.S(B2:B8)/(SQRT(COUNTA(B2:B9)))"

$ws.Range("B12").Formula = "=C2+C3+C4+C5+C6+C7+C8+C9"
$ws.Range("C12").Formula = "=AVERAGE(C2:C9)"
$ws.Range("D12").Formula = "=_xlfn.STDEV.S(C2:C8)/(SQRT(COUNTA(C2:C9)))"

$ws.Range("B13").Formula = "=D2+D3+D4+D5+D6+D7+D8+D9"
$ws.Range("C13").Formula = "=AVERAGE(D2:D9)"
$ws.Range("D13").Formula = "=_xlfn.STDEV.S(D2:D8)/(SQRT(COUNTA(D2:D9)))"

$ws.Range("B14").Formula = "=E2+E3+E4+E5+E6+E7+E8+E9"
$ws.Range("C14").Formula = "=AVERAGE(E2:E9)"
$ws.Range("D14").Formula = "=_xlfn.STDEV.S(E2:E8)/(SQRT(COUNTA(E2:E9)))"

$ws.Range("B15").Formula = "=B14*3"
$ws.Range("C15").Formula = "=AVERAGE(F2:F8)"
$ws.Range("D15").Formula = "=_xlfn.STDEV.S(F2:F8)/(SQRT(COUNTA(F2:F9)))"

# ---------------------------------------------------------------------------
# Update the sheet selection + recalc the dimension is handled automatically
# ---------------------------------------------------------------------------
$ws.Range("E30").Select()

Write-Host "edit complete"
